$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix parse-int issue: B4 and C4 were stored as text "11" (shared string),
# convert them to proper numeric values 12. D4 ("0") remains a text value
# and its shared-string index auto-shifts once the now-unused "11" string
# entry is dropped from the shared strings table.
$ws.Range("B4").Value = 12
$ws.Range("C4").Value = 12
